# "Refined metadata to be additional tab"
#
# 1. Update the time_taken timestamps (column F) on the "data" sheet.
# 2. Add a new "metadata" worksheet right after "data" with a header row
#    and a single data row describing the panel query that produced the
#    "data" sheet.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# --- 1. refresh the per-row "time_taken" timestamps on the data sheet ---
$data.Range("F2").Value = "2021-10-05 14:20:15.981474"
$data.Range("F3").Value = "2021-10-05 14:20:15.981483"
$data.Range("F4").Value = "2021-10-05 14:20:15.981486"
$data.Range("F5").Value = "2021-10-05 14:20:15.981489"
$data.Range("F6").Value = "2021-10-05 14:20:15.981492"
$data.Range("F7").Value = "2021-10-05 14:20:15.981495"

# --- 2. add the new "metadata" sheet after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (bold / bordered / centered, same style as the "data" sheet headers)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$header = $meta.Range("B1:G1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Data row
$meta.Range("A2").Value = 0
$meta.Range("A2").Font.Bold = $true
$meta.Range("A2").HorizontalAlignment = -4108
$meta.Range("A2").VerticalAlignment = -4160
$meta.Range("A2").Borders.LineStyle = 1

$meta.Range("B2").Value = "Familial hypercholesterolaemia - targeted panel"
$meta.Range("C2").Value = 772
# data_version is stored as literal text ("1.9"), not a number
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.9"
$meta.Range("E2").Value = "2021-03-23T17:56:51.816773Z"
$meta.Range("F2").Value = "2021-10-05 14:20:15.977990"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/772/?format=json"

# Restore "data" as the active sheet (workbook default view stays unchanged)
$data.Activate()
